$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Login with valid username and password", "FAILED", "chrome"),
    @("Login with valid username and password", "FAILED", "chrome"),
    @("Add New Bank Accounts", "FAILED", "chrome"),
    @("Edit The Bank Accounts", "FAILED", "chrome"),
    @("Delete The Bank Accounts", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome")
)

$startRow = 121
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
